$wb = $excel.ActiveWorkbook

# Updated "想去人数" (interested count) values for the 展览 and 全部类型 sheets.
$updates = @{
    3  = 223
    4  = 258
    6  = 231
    7  = 5738
    10 = 95
    15 = 312
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
